$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 5861
$ws.Range("H3").Value = 8348
$ws.Range("I3").Value = 6096
$ws.Range("F4").Value = 1868
$ws.Range("I4").Value = 1403
$ws.Range("I5").Value = 562
$ws.Range("I6").Value = 6888
$ws.Range("F7").Value = 24057
$ws.Range("H7").Value = 25985
$ws.Range("I7").Value = 20810

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 41
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 199
$ws.Range("I4").Value = 37
$ws.Range("I6").Value = 192
$ws.Range("I7").Value = 663

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 296
$ws.Range("I7").Value = 808

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 65
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 162
$ws.Range("I7").Value = 489

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 654
$ws.Range("I8").Value = 1253
$ws.Range("I11").Value = 311
$ws.Range("I14").Value = 117
$ws.Range("I15").Value = 233
$ws.Range("I16").Value = 60
$ws.Range("I18").Value = 153
$ws.Range("I19").Value = 573
$ws.Range("I20").Value = 514
$ws.Range("I22").Value = 55
$ws.Range("I23").Value = 208
$ws.Range("I24").Value = 59
$ws.Range("I25").Value = 109
$ws.Range("I27").Value = 188
$ws.Range("I29").Value = 1300
$ws.Range("I30").Value = 72
$ws.Range("I33").Value = 943
$ws.Range("I36").Value = 276
$ws.Range("I37").Value = 663
$ws.Range("I41").Value = 86
$ws.Range("I42").Value = 703
$ws.Range("I46").Value = 44
$ws.Range("I51").Value = 243
$ws.Range("I52").Value = 448
$ws.Range("I53").Value = 216
$ws.Range("I55").Value = 229
$ws.Range("I56").Value = 24
$ws.Range("I57").Value = 84
$ws.Range("I60").Value = 113
$ws.Range("F63").Value = 158
$ws.Range("H63").Value = 223
$ws.Range("I63").Value = 69
$ws.Range("I64").Value = 179
$ws.Range("I65").Value = 489
$ws.Range("I67").Value = 808
$ws.Range("I69").Value = 46
$ws.Range("I73").Value = 190
$ws.Range("I76").Value = 300
$ws.Range("I78").Value = 285
$ws.Range("I79").Value = 585
$ws.Range("I83").Value = 447
$ws.Range("I84").Value = 181
$ws.Range("I85").Value = 940
$ws.Range("I87").Value = 48
$ws.Range("I88").Value = 188
$ws.Range("I91").Value = 226
$ws.Range("I92").Value = 56
$ws.Range("I94").Value = 219
$ws.Range("I95").Value = 320
$ws.Range("I96").Value = 226
$ws.Range("I98").Value = 146
$ws.Range("F101").Value = 24057
$ws.Range("H101").Value = 25985
$ws.Range("I101").Value = 20810

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 161
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 447

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 116
$ws.Range("I7").Value = 320

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 357
$ws.Range("I7").Value = 943

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 380
$ws.Range("I3").Value = 449
$ws.Range("I4").Value = 66
$ws.Range("I7").Value = 1300

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 195
$ws.Range("I7").Value = 573

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 140
$ws.Range("I7").Value = 300

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 262
$ws.Range("I3").Value = 363
$ws.Range("I7").Value = 940

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 178
$ws.Range("I3").Value = 233
$ws.Range("I7").Value = 703

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 285

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I6").Value = 73
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 59

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 60
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 208

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 83
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I6").Value = 171
$ws.Range("I7").Value = 585

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 145
$ws.Range("I4").Value = 32
$ws.Range("I6").Value = 178
$ws.Range("I7").Value = 514

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 276

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 123
$ws.Range("I4").Value = 36
$ws.Range("I7").Value = 448

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 37
$ws.Range("I6").Value = 127
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 41
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 311

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 62
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 384
$ws.Range("I3").Value = 355
$ws.Range("I6").Value = 406
$ws.Range("I7").Value = 1253

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 243

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 39
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 48
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 216

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 171
$ws.Range("I7").Value = 654

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 60
